$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data values for rows 2-6 (columns A-I); column J is the shared string "train_dim2_1"
$data = @(
    @(1, 0, 3, 4, 5, 4, 2, 54, 5),
    @(2, 1, 2, 6, 3, 5, 1, 65, 5),
    @(3, 1, 4, 2, 9, 1, 5, 21, 5),
    @(4, 0, 1, 3, 4, 3, 3, 43, 5),
    @(5, 3, 4, 5, 8, 2, 4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}
